$d = $word.ActiveDocument

# Paragraph 1: "Yerba Buena, 22 de Abril de 1992"
$p1 = $d.Paragraphs(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# Paragraph 2: "ORDENANZA Nº 460"
$p2 = $d.Paragraphs(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Bold = 1

# Paragraph 3: "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
$p3 = $d.Paragraphs(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 18
$p3.Format.SpaceAfter = 18
$p3.Format.LeftIndent = 99.2
$p3.Format.RightIndent = 99.2
$p3.Range.Bold = 1

# Paragraphs 4-7: "ARTICULO PRIMERO/SEGUNDO/TERCERO/CUARTO"
$articulos = @("ARTICULO PRIMERO", "ARTICULO SEGUNDO", "ARTICULO TERCERO", "ARTICULO CUARTO")
$paraIndex = 4
foreach ($art in $articulos) {
    $p = $d.Paragraphs($paraIndex)
    $p.Format.KeepWithNext = $true
    $p.Format.SpaceAfter = 6

    $rng = $d.Content
    $found = $rng.Find.Execute($art, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.End = $rng.End + 1
    $rng.Underline = 1

    $paraIndex = $paraIndex + 1
}

Write-Output "done"
